$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Title: "Journal de bord - Sprint 0" -> "Journal de bord - Sprint 1"
#    split into two runs ("...Sprint " + "1") with identical formatting,
#    matching the target OOXML which keeps two <w:r> elements.
# ------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$titleEnd = $p1.Range.End   # position right after the paragraph mark
# last visible character of the title run (the trailing "0")
$lastChar = $d.Range($titleEnd - 2, $titleEnd - 1)
$lastChar.Text = "1"
# force a genuine run break at that position even though formatting matches
$lastChar2 = $d.Range($titleEnd - 2, $titleEnd - 1)
$lastChar2.Font.Bold = 0
$lastChar3 = $d.Range($titleEnd - 2, $titleEnd - 1)
$lastChar3.Font.Bold = 1

# ------------------------------------------------------------------
# 2) Merge the "27-03:" / long narrative / blank / "31-03:" paragraphs
#    into a single paragraph, then rewrite its digits to "10-04:".
# ------------------------------------------------------------------

# 2a. Drop the long narrative text, keep its paragraph mark for now.
$pLong = $d.Paragraphs(3)
$longTextOnly = $d.Range($pLong.Range.Start, $pLong.Range.End - 1)
$longTextOnly.Delete()

# 2b. Merge "27-03:" paragraph with the (now empty) former-long paragraph.
$pDate = $d.Paragraphs(2)
$mark1 = $d.Range($pDate.Range.End - 1, $pDate.Range.End)
$mark1.Delete()

# 2c. Merge in the blank paragraph that used to sit before "31-03:".
$pDateB = $d.Paragraphs(2)
$mark2 = $d.Range($pDateB.Range.End - 1, $pDateB.Range.End)
$mark2.Delete()

# 2d. Merge "31-03:" paragraph into the same paragraph as "27-03:".
$pDateC = $d.Paragraphs(2)
$mark3 = $d.Range($pDateC.Range.End - 1, $pDateC.Range.End)
$mark3.Delete()

# Paragraph 2 now reads "27-03:31-03:" across 4 runs: "27-03:", "31", "-03", ":"
$pMerged = $d.Paragraphs(2)
$mergedStart = $pMerged.Range.Start   # points at the <w:br/>
# Layout relative to mergedStart: +0 br, +1.. "27-03:", then "31","-03",":" then para mark
$run3 = $d.Range($mergedStart + 9, $mergedStart + 12)   # "-03" -> "4"
$run3.Text = "4"
$run2 = $d.Range($mergedStart + 7, $mergedStart + 9)    # "31" -> "-0"
$run2.Text = "-0"
$run1 = $d.Range($mergedStart + 1, $mergedStart + 7)    # "27-03:" -> "10"
$run1.Text = "10"

# Re-split the digits/punctuation back into distinct runs ("10","-0","4",":")
# using a harmless Bold on/off toggle so saved formatting stays identical.
$b1 = $d.Range($mergedStart + 3, $mergedStart + 5)   # "-0"
$b1.Font.Bold = 1
$b1 = $d.Range($mergedStart + 3, $mergedStart + 5)
$b1.Font.Bold = 0

$b2 = $d.Range($mergedStart + 5, $mergedStart + 6)   # "4"
$b2.Font.Bold = 1
$b2 = $d.Range($mergedStart + 5, $mergedStart + 6)
$b2.Font.Bold = 0

$b3 = $d.Range($mergedStart + 6, $mergedStart + 7)   # ":"
$b3.Font.Bold = 1
$b3 = $d.Range($mergedStart + 6, $mergedStart + 7)
$b3.Font.Bold = 0

# ------------------------------------------------------------------
# 3) "Ajout détails de Jira au résumé Scrum de l'équipe"
#    -> "Production de documents Product Owner."
# ------------------------------------------------------------------
$rFind = $d.Content
$rFind.Find.ClearFormatting()
$rFind.Find.Execute("Ajout détails de Jira au résumé Scrum de l'équipe", $true, $false, $false, $false, $false, `
                     $true, 1, $false, "Production de documents Product Owner.", 2)

# ------------------------------------------------------------------
# 4) Drop "Ajout détails aux normes" and "Résumé des vidéos effectués",
#    leaving a single empty trailing paragraph.
# ------------------------------------------------------------------
$idxNormes = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*normes*") { $idxNormes = $i }
}
$pNormes = $d.Paragraphs($idxNormes)
$markNormes = $d.Range($pNormes.Range.End - 1, $pNormes.Range.End)
$markNormes.Delete()

$pTail = $d.Paragraphs($idxNormes)
$tailText = $d.Range($pTail.Range.Start, $pTail.Range.End - 1)
$tailText.Delete()
